$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.009.57"
$ws.Range("E2").Value = "  -1.77%  "

$ws.Range("D3").Value = "2.300.63"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.05"
$ws.Range("E5").Value = "  -1.12%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "104.26"
$ws.Range("E6").Value = "  -0.69%  "

$ws.Range("E7").Value = "  -1.50%  "

$ws.Range("E8").Value = "  -0.13%  "

$ws.Range("E9").Value = "  -1.19%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.85"
$ws.Range("E10").Value = "  -3.21%  "

$ws.Range("E11").Value = "  -1.13%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  +0.54%  "

$ws.Range("E13").Value = "  +1.12%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.976"
$ws.Range("E14").Value = "  -2.51%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.41"
$ws.Range("E15").Value = "  -3.75%  "

$ws.Range("D16").Value = "2.648.88"
$ws.Range("E16").Value = "  -2.07%  "

$ws.Range("D17").Value = "2.293.10"
$ws.Range("E17").Value = "  -2.35%  "

$ws.Range("D18").Value = "42.016.18"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.72"
$ws.Range("E19").Value = "  +0.42%  "

$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "73.15"
$ws.Range("E21").Value = "  -5.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.59"
$ws.Range("E22").Value = "  -1.10%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "258.97"
$ws.Range("E23").Value = "  -0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.31"
$ws.Range("E24").Value = "  -0.52%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.94"
$ws.Range("E25").Value = "  +3.41%  "

$ws.Range("E26").Value = "  +0.72%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.95"
$ws.Range("E27").Value = "  -3.70%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.88"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("E29").Value = "  +0.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.15"
$ws.Range("E30").Value = "  +0.00%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.07"
$ws.Range("E31").Value = "  -6.09%  "

$ws.Range("E32").Value = "  -0.05%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.91"
$ws.Range("E33").Value = "  -2.86%  "

$ws.Range("E34").Value = "  -4.42%  "

$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.120"
$ws.Range("E35").Value = "  +3.56%  "

$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.130"
$ws.Range("E36").Value = "  -0.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.64"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("E38").Value = "  +9.64%  "

$ws.Range("E39").Value = "  -2.48%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.69"
$ws.Range("E40").Value = "  -2.84%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.18"
$ws.Range("E41").Value = "  +15.47%  "

$ws.Range("E42").Value = "  +0.94%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "71.07"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("E44").Value = "  -1.65%  "

$ws.Range("E45").Value = "  -0.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.17"
$ws.Range("E46").Value = "  +2.23%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "114.75"
$ws.Range("E47").Value = "  -0.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "79.07"
$ws.Range("E48").Value = "  +6.37%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.06"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("E50").Value = "  -3.25%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.29"
$ws.Range("E51").Value = "  +2.13%  "

